# Fruta / hortaliza, semanal
# The underlying data rows (2-26) get reshuffled into a new order (the
# upstream weekly refresh reassigns each record to a different row), while
# the column layout (D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, O=Origen, P=Precio $/Kg) stays the same.
#
# Below, for every destination row we list the final values that row must
# contain after the refresh (these are simply the previous rows' values,
# redistributed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colD = 4   # Fecha
$colJ = 10  # Volumen
$colK = 11  # Precio minimo
$colL = 12  # Precio maximo
$colM = 13  # Precio promedio ponderado
$colO = 15  # Origen
$colP = 16  # Precio $/Kg

# row, Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$data = @(
    @(2,  44461, 100, 13000, 14000, 13500, "Provincia del Elquí", 540),
    @(3,  44537, 160, 8500,  9000,  8719,  "Región del Maule", 349),
    @(4,  44482, 430, 8000,  8500,  8267,  "Región de O'Higgins", 331),
    @(5,  44505, 180, 6000,  6500,  6222,  "Región del Maule", 249),
    @(6,  44384, 100, 12000, 13000, 12500, "Región de Coquimbo", 500),
    @(7,  44503, 250, 9000,  10000, 9400,  "Provincia de Melipilla", 376),
    @(8,  44509, 100, 6500,  7000,  6750,  "Región Metropolitana", 270),
    @(9,  44526, 100, 7500,  8000,  7750,  "Región Metropolitana", 310),
    @(10, 44545, 140, 14000, 15000, 14429, "Provincia de Chacabuco", 577),
    @(11, 44540, 140, 11000, 12000, 11429, "Región del Maule", 457),
    @(12, 44454, 100, 13000, 14000, 13500, "Provincia del Elquí", 540),
    @(13, 44497, 150, 6000,  6500,  6333,  "Región Metropolitana", 253),
    @(14, 44351, 100, 15000, 16000, 15500, "Región Metropolitana", 620),
    @(15, 44476, 100, 7000,  7500,  7250,  "Región Metropolitana", 290),
    @(16, 44316, 100, 16000, 18000, 17000, "Región Metropolitana", 680),
    @(17, 44188, 100, 18000, 20000, 19000, "Región Metropolitana", 760),
    @(18, 44523, 100, 9000,  10000, 9500,  "Región Metropolitana", 380),
    @(19, 44483, 350, 5500,  6000,  5714,  "Región Metropolitana", 229),
    @(20, 44533, 180, 8000,  8500,  8222,  "Región del Maule", 329),
    @(21, 44160, 100, 9000,  10000, 9500,  "Región Metropolitana", 380),
    @(22, 44162, 100, 7500,  8000,  7750,  "Región Metropolitana", 310),
    @(23, 44498, 220, 7000,  7500,  7273,  "Región Metropolitana", 291),
    @(24, 44517, 130, 6000,  6500,  6269,  "Región Metropolitana", 251),
    @(25, 44335, 100, 18000, 20000, 19000, "Provincia de Limarí", 760),
    @(26, 44467, 100, 8000,  9000,  8500,  "Región Metropolitana", 340)
)

foreach ($row in $data) {
    $r      = $row[0]
    $fecha  = $row[1]
    $vol    = $row[2]
    $pmin   = $row[3]
    $pmax   = $row[4]
    $pprom  = $row[5]
    $origen = $row[6]
    $pkg    = $row[7]

    $ws.Cells.Item($r, $colD).Value2 = $fecha
    $ws.Cells.Item($r, $colJ).Value2 = $vol
    $ws.Cells.Item($r, $colK).Value2 = $pmin
    $ws.Cells.Item($r, $colL).Value2 = $pmax
    $ws.Cells.Item($r, $colM).Value2 = $pprom
    $ws.Cells.Item($r, $colO).Value2 = $origen
    $ws.Cells.Item($r, $colP).Value2 = $pkg
}
